# Trade #89 closed at 2026-02-17 21:18:37 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook to reflect:
#   - Trade #117 (row 118 in "All Trades" / row 85 in "MarketMaking") being
#     closed (early_exit), with its exit price / P&L filled in.
#   - A brand-new open Trade #150 appended to both the "All Trades" sheet
#     (as row 151) and the per-strategy "MarketMaking" sheet (as row 118).
#   - Updated roll-up metrics on the "Summary" and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.06   # Current Capital
$wsSummary.Range("B4").Value = 0.85      # Total P&L $
$wsSummary.Range("B6").Value = 117       # Total Trades
$wsSummary.Range("B8").Value = 45        # Losing Trades
$wsSummary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.06     # Capital
$wsStatus.Range("D5").Value = 84         # Trades
$wsStatus.Range("E5").Value = 0.74       # P&L $
$wsStatus.Range("F5").Value = 1.06       # P&L %
$wsStatus.Range("G5").Value = 45.24      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Trade #117 (row 118) transitions from OPEN to CLOSED (early_exit)
$wsAll.Cells.Item(118, 7).Value = 0.91            # Exit Price
$wsAll.Cells.Item(118, 8).Value = "CLOSED"        # Status
$wsAll.Cells.Item(118, 9).Value = -2.0761         # P&L %
$wsAll.Cells.Item(118, 10).Value = -0.02          # P&L $
$wsAll.Cells.Item(118, 11).Value = 101.06         # Capital After
$wsAll.Cells.Item(118, 12).Value = "early_exit"   # Exit Reason
$wsAll.Cells.Item(118, 13).Value = 0.14           # Duration (min)

# New row 151: Trade #150 (freshly opened)
$wsAll.Cells.Item(151, 1).Value = 150
$wsAll.Cells.Item(151, 2).NumberFormat = "@"
$wsAll.Cells.Item(151, 2).Value = "2026-02-17"
$wsAll.Cells.Item(151, 2).Style = "Normal"
$wsAll.Cells.Item(151, 3).Value = "21:18:31"
$wsAll.Cells.Item(151, 4).Value = "MarketMaking"
$wsAll.Cells.Item(151, 5).Value = "DOWN"
$wsAll.Cells.Item(151, 6).Value = 0.929293
$wsAll.Cells.Item(151, 8).Value = "OPEN"
$wsAll.Cells.Item(151, 9).Value = 0
$wsAll.Cells.Item(151, 10).Value = 0
$wsAll.Cells.Item(151, 11).Value = 101.0796151053151
$wsAll.Cells.Item(151, 13).Value = 0
$wsAll.Cells.Item(151, 14).Value = 0
$wsAll.Cells.Item(151, 15).Value = 0
$wsAll.Cells.Item(151, 16).Value = 0.6
$wsAll.Cells.Item(151, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet (per-strategy trade log; different column layout)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Trade #117 (row 85) transitions from OPEN to CLOSED (early_exit)
$wsMM.Cells.Item(85, 7).Value = 0.91              # Exit Price
$wsMM.Cells.Item(85, 8).Value = "CLOSED"          # Status
$wsMM.Cells.Item(85, 9).Value = -2.0761           # P&L %
$wsMM.Cells.Item(85, 10).Value = -0.02            # P&L $
$wsMM.Cells.Item(85, 11).Value = 101.06           # Capital After
$wsMM.Cells.Item(85, 16).Value = "early_exit"     # Exit Reason
$wsMM.Cells.Item(85, 17).Value = 0.14             # Duration (min)

# New row 118: Trade #150 (freshly opened)
$wsMM.Cells.Item(118, 1).Value = 150
$wsMM.Cells.Item(118, 2).NumberFormat = "@"
$wsMM.Cells.Item(118, 2).Value = "2026-02-17"
$wsMM.Cells.Item(118, 2).Style = "Normal"
$wsMM.Cells.Item(118, 3).Value = "21:18:31"
$wsMM.Cells.Item(118, 4).Value = "MarketMaking"
$wsMM.Cells.Item(118, 5).Value = "DOWN"
$wsMM.Cells.Item(118, 6).Value = 0.929293
$wsMM.Cells.Item(118, 8).Value = "OPEN"
$wsMM.Cells.Item(118, 9).Value = 0
$wsMM.Cells.Item(118, 10).Value = 0
$wsMM.Cells.Item(118, 11).Value = 101.0796151053151
$wsMM.Cells.Item(118, 12).Value = 0
$wsMM.Cells.Item(118, 13).Value = 0
$wsMM.Cells.Item(118, 14).Value = 0.6
$wsMM.Cells.Item(118, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(118, 17).Value = 0
